# Initial implementation of 'change password' functionality
# Adds two new error-code rows (10034, 10035) to the error codes sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36 -> error code 10034 : user viewing own record (Information)
$ws.Range("B36").Value = "message_10034_user_viewing_own_record"
$ws.Range("D36").Value = "Information"

# Row 37 -> error code 10035 : password not long enough (Error)
$ws.Range("B37").Value = "message_10035_password_not_long_enough"
$ws.Range("D37").Value = "Error"

# Update the remembered selection on the sheet (matches author's last position)
$ws.Range("D38").Select()
